# Auto-generated edit script applying numeric corrections per commit diff
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 45
$ws.Range("I9").Value = 45
$ws.Range("K9").Value = 45
$ws.Range("M9").Value = 124
$ws.Range("H12").Value = 383.14285
$ws.Range("J12").Value = 310.16666
$ws.Range("L12").Value = 310.16666
$ws.Range("N12").Value = -650.16666
$ws.Range("H70").Value = 682912.0600000001
$ws.Range("I70").Value = 1136064.5
$ws.Range("K70").Value = 3408193.5
$ws.Range("M70").Value = -3407923.5
$ws.Range("H73").Value = 682912.0600000001
$ws.Range("I73").Value = 1136064.5
$ws.Range("K73").Value = 3408193.5
$ws.Range("M73").Value = -3407257.5
$ws.Range("H76").Value = 9375727
$ws.Range("I76").Value = 6067.6665
$ws.Range("J76").Value = 14060557
$ws.Range("K76").Value = 6067.6665
$ws.Range("L76").Value = 14060557
$ws.Range("M76").Value = -5752.6665
$ws.Range("N76").Value = -14061187
$ws.Range("H79").Value = 9375727
$ws.Range("I79").Value = 6067.6665
$ws.Range("J79").Value = 14060557
$ws.Range("K79").Value = 6067.6665
$ws.Range("L79").Value = 14060557
$ws.Range("M79").Value = -4975.6665
$ws.Range("N79").Value = -14062741
$ws.Range("H86").Value = 9574949
$ws.Range("I86").Value = 1783.5
$ws.Range("J86").Value = 18277826
$ws.Range("K86").Value = 1783.5
$ws.Range("L86").Value = 18277826
$ws.Range("M86").Value = -660.5
$ws.Range("N86").Value = -18280072
$ws.Range("H89").Value = 9574949
$ws.Range("I89").Value = 1783.5
$ws.Range("J89").Value = 18277826
$ws.Range("K89").Value = 8917.5
$ws.Range("L89").Value = 91389130
$ws.Range("M89").Value = -3301.5
$ws.Range("N89").Value = -91400362
$ws.Range("H100").Value = 4180.778
$ws.Range("I100").Value = 4682.3076
$ws.Range("J100").Value = 2876.8
$ws.Range("K100").Value = 4682.3076
$ws.Range("L100").Value = 2876.8
$ws.Range("M100").Value = -4141.3076
$ws.Range("N100").Value = -3958.8
$ws.Range("H103").Value = 569.8
$ws.Range("J103").Value = 699.6667
$ws.Range("L103").Value = 2099.0001
$ws.Range("N103").Value = -3271.0001
$ws.Range("H111").Value = 28559.125
$ws.Range("I111").Value = 1439.4546
$ws.Range("K111").Value = 4318.3638
$ws.Range("M111").Value = -1251.3638
$ws.Range("H138").Value = 1733.5962
$ws.Range("I138").Value = 1221.5883
$ws.Range("J138").Value = 1982.2858
$ws.Range("K138").Value = 3664.7649
$ws.Range("L138").Value = 5946.857400000001
$ws.Range("M138").Value = 1475.2351
$ws.Range("N138").Value = -16226.8574
$ws.Range("H141").Value = 4900
$ws.Range("I141").Value = 5348.75
$ws.Range("J141").Value = 4002.5
$ws.Range("K141").Value = 16046.25
$ws.Range("L141").Value = 12007.5
$ws.Range("M141").Value = -10866.25
$ws.Range("N141").Value = -22367.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3191.9644
$ws.Range("I32").Value = 2606.463
$ws.Range("J32").Value = 19000.5
$ws.Range("K32").Value = 2606.463
$ws.Range("L32").Value = 19000.5
$ws.Range("M32").Value = -2319.463
$ws.Range("N32").Value = -19574.5
$ws.Range("H45").Value = 4249.6665
$ws.Range("I45").Value = 4199.5
$ws.Range("K45").Value = 4199.5
$ws.Range("M45").Value = -3822.5
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()
$ws.Range("H61").Value = 2283.75
$ws.Range("I61").Value = 2009.7222
$ws.Range("K61").Value = 2009.7222
$ws.Range("M61").Value = -1797.7222
$ws.Range("H74").Value = 27781850
$ws.Range("I74").Value = 33337144
$ws.Range("K74").Value = 33337144
$ws.Range("M74").Value = -33336270
$ws.Range("H77").Value = 27781850
$ws.Range("I77").Value = 33337144
$ws.Range("K77").Value = 166685720
$ws.Range("M77").Value = -166681352
$ws.Range("H81").Value = 49999.5
$ws.Range("I81").Value = 49999.5
$ws.Range("K81").Value = 49999.5
$ws.Range("M81").Value = -49001.5
$ws.Range("H84").Value = 49999.5
$ws.Range("I84").Value = 49999.5
$ws.Range("K84").Value = 149998.5
$ws.Range("M84").Value = -145006.5
$ws.Range("H110").Value = 22728164
$ws.Range("I110").Value = 29412430
$ws.Range("K110").Value = 29412430
$ws.Range("M110").Value = -29410385
$ws.Range("H132").Value = 3508.9375
$ws.Range("I132").Value = 1741.5454
$ws.Range("K132").Value = 5224.6362
$ws.Range("M132").Value = -2694.6362
$ws.Range("H136").Value = 2283.75
$ws.Range("I136").Value = 2009.7222
$ws.Range("K136").Value = 6029.1666
$ws.Range("M136").Value = -3479.1666

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 260.14285
$ws.Range("I22").Value = 253.5
$ws.Range("K22").Value = 253.5
$ws.Range("M22").Value = -80.5
$ws.Range("H105").Value = 1773.4138
$ws.Range("I105").Value = 1506.15
$ws.Range("K105").Value = 1506.15
$ws.Range("M105").Value = 240.8499999999999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 821.2308
$ws.Range("I16").Value = 698.4545000000001
$ws.Range("K16").Value = 698.4545000000001
$ws.Range("M16").Value = -411.4545000000001
$ws.Range("H22").Value = 516.3333
$ws.Range("I22").Value = 499.5
$ws.Range("J22").Value = 550
$ws.Range("K22").Value = 499.5
$ws.Range("L22").Value = 550
$ws.Range("M22").Value = -149.5
$ws.Range("N22").Value = -1250
$ws.Range("H31").Value = 3193.9558
$ws.Range("I31").Value = 2050.8
$ws.Range("K31").Value = 2050.8
$ws.Range("M31").Value = -1755.8
$ws.Range("H34").Value = 3193.9558
$ws.Range("I34").Value = 2050.8
$ws.Range("K34").Value = 2050.8
$ws.Range("M34").Value = -1848.8
$ws.Range("H50").Value = 74500
$ws.Range("J50").Value = 74500
$ws.Range("L50").Value = 74500
$ws.Range("N50").Value = -75750
$ws.Range("H62").Value = 112590
$ws.Range("J62").Value = 139487.5
$ws.Range("L62").Value = 139487.5
$ws.Range("N62").Value = -140735.5
$ws.Range("H65").Value = 112590
$ws.Range("J65").Value = 139487.5
$ws.Range("L65").Value = 697437.5
$ws.Range("N65").Value = -703677.5
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("H75").Value = 35000
$ws.Range("J75").Value = 35000
$ws.Range("L75").Value = 35000
$ws.Range("N75").Value = -36996
$ws.Range("H78").Value = 35000
$ws.Range("J78").Value = 35000
$ws.Range("L78").Value = 105000
$ws.Range("N78").Value = -114984
$ws.Range("H103").Value = 16490.166
$ws.Range("I103").Value = 16490.166
$ws.Range("K103").Value = 16490.166
$ws.Range("M103").Value = -15318.166
$ws.Range("H105").Value = 481.125
$ws.Range("I105").Value = 481.125
$ws.Range("K105").Value = 481.125
$ws.Range("M105").Value = 1265.875
$ws.Range("H107").Value = 471.85
$ws.Range("I107").Value = 448.9375
$ws.Range("K107").Value = 448.9375
$ws.Range("M107").Value = 1471.0625
$ws.Range("H113").Value = 821.2308
$ws.Range("I113").Value = 698.4545000000001
$ws.Range("K113").Value = 698.4545000000001
$ws.Range("M113").Value = 1471.5455

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 6290.4443
$ws.Range("J39").Value = 4191.8667
$ws.Range("L39").Value = 12575.6001
$ws.Range("N39").Value = -13163.6001
$ws.Range("H92").Value = 698
$ws.Range("I92").Value = 698
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 2094
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = -846
$ws.Range("N92").ClearContents()
$ws.Range("H97").Value = 761
$ws.Range("I97").Value = 796.6667
$ws.Range("J97").Value = 654
$ws.Range("K97").Value = 2390.0001
$ws.Range("L97").Value = 1962
$ws.Range("M97").Value = -1894.0001
$ws.Range("N97").Value = -2954
$ws.Range("H115").Value = 1900
$ws.Range("I115").Value = 1900
$ws.Range("K115").Value = 5700
$ws.Range("M115").Value = -4525

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 30000
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 30000
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 30000
$ws.Range("M44").ClearContents()
$ws.Range("N44").Value = -31192
$ws.Range("H62").Value = 39399.668
$ws.Range("I62").Value = 39099.5
$ws.Range("K62").Value = 39099.5
$ws.Range("M62").Value = -38413.5
$ws.Range("H65").Value = 39399.668
$ws.Range("I65").Value = 39099.5
$ws.Range("K65").Value = 117298.5
$ws.Range("M65").Value = -113866.5
$ws.Range("H70").Value = 11248.75
$ws.Range("I70").Value = 11081.833
$ws.Range("K70").Value = 11081.833
$ws.Range("M70").Value = -10811.833
$ws.Range("H73").Value = 11248.75
$ws.Range("I73").Value = 11081.833
$ws.Range("K73").Value = 11081.833
$ws.Range("M73").Value = -10145.833
$ws.Range("H92").Value = 10941.167
$ws.Range("J92").Value = 10941.167
$ws.Range("L92").Value = 10941.167
$ws.Range("N92").Value = -14685.167
$ws.Range("H102").Value = 1967.1177
$ws.Range("I102").Value = 958.1429000000001
$ws.Range("K102").Value = 958.1429000000001
$ws.Range("M102").Value = 663.8570999999999
$ws.Range("H113").Value = 3037.4375
$ws.Range("I113").Value = 1085.5714
$ws.Range("K113").Value = 1085.5714
$ws.Range("M113").Value = 1084.4286

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H29").Value = 15999.5
$ws.Range("I29").Value = 15999.5
$ws.Range("K29").Value = 15999.5
$ws.Range("M29").Value = -15704.5
$ws.Range("H40").Value = 3421.4
$ws.Range("I40").Value = 3421.4
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 3421.4
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -3285.4
$ws.Range("N40").ClearContents()
$ws.Range("H46").Value = 1984.3462
$ws.Range("J46").Value = 2313.9524
$ws.Range("L46").Value = 2313.9524
$ws.Range("N46").Value = -2689.9524
$ws.Range("H93").Value = 25645256
$ws.Range("J93").Value = 4998.25
$ws.Range("L93").Value = 4998.25
$ws.Range("N93").Value = -7494.25
$ws.Range("H104").Value = 12873.8
$ws.Range("J104").Value = 12873.8
$ws.Range("L104").Value = 12873.8
$ws.Range("N104").Value = -19861.8
$ws.Range("H122").Value = 3729.9412
$ws.Range("I122").Value = 3225.5557
$ws.Range("K122").Value = 9676.667099999999
$ws.Range("M122").Value = -7226.667099999999
$ws.Range("H132").Value = 47626744
$ws.Range("I132").Value = 71432160
$ws.Range("J132").Value = 15912
$ws.Range("K132").Value = 214296480
$ws.Range("L132").Value = 47736
$ws.Range("M132").Value = -214293950
$ws.Range("N132").Value = -52796
$ws.Range("H136").Value = 4029.1
$ws.Range("I136").Value = 4422.643
$ws.Range("J136").Value = 3110.8333
$ws.Range("K136").Value = 13267.929
$ws.Range("L136").Value = 9332.499899999999
$ws.Range("M136").Value = -10717.929
$ws.Range("N136").Value = -14432.4999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H53").Value = 31000
$ws.Range("J53").Value = 31000
$ws.Range("L53").Value = 31000
$ws.Range("N53").Value = -32214
$ws.Range("H62").Value = 1197019.4
$ws.Range("J62").Value = 7116.4375
$ws.Range("L62").Value = 7116.4375
$ws.Range("N62").Value = -8364.4375
$ws.Range("H64").Value = 51665.332
$ws.Range("I64").Value = 49998.5
$ws.Range("K64").Value = 49998.5
$ws.Range("M64").Value = -49750.5
$ws.Range("H65").Value = 1197019.4
$ws.Range("J65").Value = 7116.4375
$ws.Range("L65").Value = 35582.1875
$ws.Range("N65").Value = -41822.1875
$ws.Range("H67").Value = 51665.332
$ws.Range("I67").Value = 49998.5
$ws.Range("K67").Value = 49998.5
$ws.Range("M67").Value = -49140.5
$ws.Range("H81").Value = 14290066
$ws.Range("I81").Value = 2865.125
$ws.Range("J81").Value = 33339668
$ws.Range("K81").Value = 5730.25
$ws.Range("L81").Value = 66679336
$ws.Range("M81").Value = -4669.25
$ws.Range("N81").Value = -66681458
$ws.Range("H84").Value = 14290066
$ws.Range("I84").Value = 2865.125
$ws.Range("J84").Value = 33339668
$ws.Range("K84").Value = 28651.25
$ws.Range("L84").Value = 333396680
$ws.Range("M84").Value = -23347.25
$ws.Range("N84").Value = -333407288
$ws.Range("H107").Value = 1857.5714
$ws.Range("I107").Value = 1649
$ws.Range("J107").Value = 1941
$ws.Range("K107").Value = 4947
$ws.Range("L107").Value = 5823
$ws.Range("M107").Value = -3027
$ws.Range("N107").Value = -9663
$ws.Range("H122").Value = 1731.925
$ws.Range("J122").Value = 2155.5
$ws.Range("L122").Value = 6466.5
$ws.Range("N122").Value = -11366.5
$ws.Range("H126").Value = 2272.8462
$ws.Range("I126").Value = 1522
$ws.Range("K126").Value = 4566
$ws.Range("M126").Value = -2096
$ws.Range("H136").Value = 5785.857
$ws.Range("I136").Value = 5758.0835
$ws.Range("K136").Value = 17274.2505
$ws.Range("M136").Value = -14724.2505

Write-Host "Applied all cell updates"